$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1586.4348
$ws.Range("I137").Value = 1500.2632
$ws.Range("J137").Value = 1995.75
$ws.Range("K137").Value = 4500.7896
$ws.Range("L137").Value = 5987.25
$ws.Range("M137").Value = -1950.7896
$ws.Range("N137").Value = -11087.25
$ws.Range("H138").Value = 3001.5256
$ws.Range("I138").Value = 2302.3684
$ws.Range("J138").Value = 3226.678
$ws.Range("K138").Value = 6907.1052
$ws.Range("L138").Value = 9680.034
$ws.Range("M138").Value = -1767.1052
$ws.Range("N138").Value = -19960.034

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10017.54
$ws.Range("I32").Value = 8744.632
$ws.Range("J32").Value = 34202.8
$ws.Range("K32").Value = 8744.632
$ws.Range("L32").Value = 34202.8
$ws.Range("M32").Value = -8457.632
$ws.Range("N32").Value = -34776.8
$ws.Range("H74").Value = 720.8444
$ws.Range("I74").Value = 620.3
$ws.Range("J74").Value = 1525.2
$ws.Range("K74").Value = 620.3
$ws.Range("L74").Value = 1525.2
$ws.Range("M74").Value = 253.7
$ws.Range("N74").Value = -3273.2
$ws.Range("H77").Value = 720.8444
$ws.Range("I77").Value = 620.3
$ws.Range("J77").Value = 1525.2
$ws.Range("K77").Value = 3101.5
$ws.Range("L77").Value = 7626
$ws.Range("M77").Value = 1266.5
$ws.Range("N77").Value = -16362
$ws.Range("H122").Value = 2008.4117
$ws.Range("I122").Value = 1657.5333
$ws.Range("K122").Value = 4972.5999
$ws.Range("M122").Value = -2522.5999
$ws.Range("H132").Value = 13336.313
$ws.Range("I132").Value = 16028.171
$ws.Range("J132").Value = 2299.7
$ws.Range("K132").Value = 48084.513
$ws.Range("L132").Value = 6899.099999999999
$ws.Range("M132").Value = -45554.513
$ws.Range("N132").Value = -11959.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 62739.055
$ws.Range("I86").Value = 74993.664
$ws.Range("J86").Value = 1466
$ws.Range("K86").Value = 74993.664
$ws.Range("L86").Value = 1466
$ws.Range("M86").Value = -73870.664
$ws.Range("N86").Value = -3712
$ws.Range("H89").Value = 62739.055
$ws.Range("I89").Value = 74993.664
$ws.Range("J89").Value = 1466
$ws.Range("K89").Value = 374968.32
$ws.Range("L89").Value = 7330
$ws.Range("M89").Value = -369352.32
$ws.Range("N89").Value = -18562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30310.51
$ws.Range("I31").Value = 1144.0769
$ws.Range("J31").Value = 60643.6
$ws.Range("K31").Value = 1144.0769
$ws.Range("L31").Value = 60643.6
$ws.Range("M31").Value = -849.0769
$ws.Range("N31").Value = -61233.6
$ws.Range("H34").Value = 30310.51
$ws.Range("I34").Value = 1144.0769
$ws.Range("J34").Value = 60643.6
$ws.Range("K34").Value = 1144.0769
$ws.Range("L34").Value = 60643.6
$ws.Range("M34").Value = -942.0769
$ws.Range("N34").Value = -61047.6
$ws.Range("H58").Value = 1519.4054
$ws.Range("I58").Value = 1262.4667
$ws.Range("J58").Value = 2620.5715
$ws.Range("K58").Value = 1262.4667
$ws.Range("L58").Value = 2620.5715
$ws.Range("M58").Value = -1059.4667
$ws.Range("N58").Value = -3026.5715
$ws.Range("H134").Value = 1471.2
$ws.Range("I134").Value = 856.5833
$ws.Range("J134").Value = 3929.6667
$ws.Range("K134").Value = 2569.7499
$ws.Range("L134").Value = 11789.0001
$ws.Range("M134").Value = -34.7498999999998
$ws.Range("N134").Value = -16859.0001
$ws.Range("H136").Value = 1519.4054
$ws.Range("I136").Value = 1262.4667
$ws.Range("J136").Value = 2620.5715
$ws.Range("K136").Value = 3787.4001
$ws.Range("L136").Value = 7861.7145
$ws.Range("M136").Value = -1237.4001
$ws.Range("N136").Value = -12961.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 638.625
$ws.Range("I2").Value = 16.285715
$ws.Range("J2").Value = 1122.6666
$ws.Range("K2").Value = 97.71429
$ws.Range("L2").Value = 6735.9996
$ws.Range("M2").Value = 15.28570999999999
$ws.Range("N2").Value = -6961.9996
$ws.Range("H23").Value = 486.66666
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 486.66666
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 1459.99998
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -1929.99998
$ws.Range("H33").Value = 1636.7894
$ws.Range("I33").Value = 1707
$ws.Range("J33").Value = 1604.3846
$ws.Range("K33").Value = 10242
$ws.Range("L33").Value = 9626.3076
$ws.Range("M33").Value = -9959
$ws.Range("N33").Value = -10192.3076
$ws.Range("H56").Value = 4247.4
$ws.Range("I56").Value = 4247.4
$ws.Range("K56").Value = 4247.4
$ws.Range("M56").Value = -3717.4
$ws.Range("H131").Value = 734365.4
$ws.Range("I131").Value = 686.36365
$ws.Range("J131").Value = 873511.4
$ws.Range("K131").Value = 2059.09095
$ws.Range("L131").Value = 2620534.2
$ws.Range("M131").Value = 2980.90905
$ws.Range("N131").Value = -2630614.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1969.2727
$ws.Range("I43").Value = 998.0741
$ws.Range("J43").Value = 6339.6665
$ws.Range("K43").Value = 998.0741
$ws.Range("L43").Value = 6339.6665
$ws.Range("M43").Value = -847.0741
$ws.Range("N43").Value = -6641.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9717701
$ws.Range("I16").Value = 14003424
$ws.Range("J16").Value = 2003400
$ws.Range("K16").Value = 14003424
$ws.Range("L16").Value = 2003400
$ws.Range("M16").Value = -14003254
$ws.Range("N16").Value = -2003740
$ws.Range("H40").Value = 61996.41
$ws.Range("I40").Value = 253624.75
$ws.Range("J40").Value = 3033.8462
$ws.Range("K40").Value = 253624.75
$ws.Range("L40").Value = 3033.8462
$ws.Range("M40").Value = -253488.75
$ws.Range("N40").Value = -3305.8462
$ws.Range("H55").Value = 307957.16
$ws.Range("I55").Value = 875553.6
$ws.Range("J55").Value = 509.08334
$ws.Range("K55").Value = 875553.6
$ws.Range("L55").Value = 509.08334
$ws.Range("M55").Value = -875380.6
$ws.Range("N55").Value = -855.08334
$ws.Range("H132").Value = 4815.143
$ws.Range("I132").Value = 9347
$ws.Range("J132").Value = 3002.4
$ws.Range("K132").Value = 28041
$ws.Range("L132").Value = 9007.2
$ws.Range("M132").Value = -25511
$ws.Range("N132").Value = -14067.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7759.5527
$ws.Range("I132").Value = 5740.645
$ws.Range("J132").Value = 16700.428
$ws.Range("K132").Value = 17221.935
$ws.Range("L132").Value = 50101.284
$ws.Range("M132").Value = -14691.935
$ws.Range("N132").Value = -55161.284
